$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new row at row 265 (pushes old row265.. down by one)
$ws.Rows("265:265").Insert()

# Insert second new row at row 273 (pushes old row272-now-at-273 .. down by one more)
$ws.Rows("273:273").Insert()

# Fill in the brand-new row 265 with its full record
$ws.Range("A265").Value = 10
$ws.Range("B265").Value = "Vega Modelo de Temuco"
$ws.Range("C265").Value = "La Araucanía"
$ws.Range("D265").Value = 44846
$ws.Range("E265").Value = 9
$ws.Range("F265").Value = 100112017
$ws.Range("G265").Value = "Apio"
$ws.Range("H265").Value = "Americana (o)"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 125
$ws.Range("K265").Value = 10000
$ws.Range("L265").Value = 10000
$ws.Range("M265").Value = 10000
$ws.Range("N265").Value = "$/docena de matas"
$ws.Range("O265").Value = "Provincia del Elquí"
$ws.Range("P265").Value = 1667
$ws.Range("Q265").Value = 6
$ws.Range("R265").Value = "Hortaliza"

# Fill in the brand-new row 273 with its full record
$ws.Range("A273").Value = 10
$ws.Range("B273").Value = "Vega Modelo de Temuco"
$ws.Range("C273").Value = "La Araucanía"
$ws.Range("D273").Value = 44845
$ws.Range("E273").Value = 9
$ws.Range("F273").Value = 100112017
$ws.Range("G273").Value = "Apio"
$ws.Range("H273").Value = "Americana (o)"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 120
$ws.Range("K273").Value = 9000
$ws.Range("L273").Value = 10000
$ws.Range("M273").Value = 9458
$ws.Range("N273").Value = "$/docena de matas"
$ws.Range("O273").Value = "Provincia del Elquí"
$ws.Range("P273").Value = 1576
$ws.Range("Q273").Value = 6
$ws.Range("R273").Value = "Hortaliza"
